$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Simple RMTL -> PMTL rename on the sheets that only show the column
#    header (display name) row and do not get their columns reordered.
# ---------------------------------------------------------------------------
$simpleRenames = @(
    @{ Sheet = "SNV gene-level";               Cell = "D1" },
    @{ Sheet = "CNV gene-level";                Cell = "P1" },
    @{ Sheet = "Fusion gene-level";             Cell = "G1" },
    @{ Sheet = "Fusion fusion-level";           Cell = "T1" },
    @{ Sheet = "TPM stats gene-wise z-scores";  Cell = "D1" },
    @{ Sheet = "TPM stats group-wise z-scores"; Cell = "D1" }
)

foreach ($item in $simpleRenames) {
    $ws = $wb.Worksheets.Item($item.Sheet)
    $ws.Range($item.Cell).Value = "PMTL"
}

# ---------------------------------------------------------------------------
# 2) "SNV variant-level" sheet: move the "Variant ID hg38" / "Variant_ID_hg38"
#    column (originally column I) and the "Protein change" / "Protein_change"
#    column (originally column T) so that they sit right after the
#    "Gene symbol" / "Gene_symbol" column (column C), becoming the new
#    columns D and E. Everything that used to sit between D and S shifts two
#    slots to the right.
#
#    Using whole-column Cut + Insert (instead of writing values back through
#    Value/Value2) keeps literal text such as "0.19%" as text instead of
#    Excel re-interpreting it as a percentage number, and keeps the original
#    (style-less) cell formatting intact.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SNV variant-level")

# Move "Variant ID hg38" (column I) to just after "Gene symbol" (-> column D).
$ws2.Range("I:I").Cut()
$ws2.Range("D:D").Insert()

# "Protein change" is still in column T (its absolute position didn't move:
# the column removed from I and the column inserted before D are both to the
# left of T, so the shifts cancel out). Move it to just after the now-placed
# "Variant ID hg38" (-> column E).
$ws2.Range("T:T").Cut()
$ws2.Range("E:E").Insert()

# Row 1 display-name cell for the RMTL/PMTL marker now lives in column F
# (the JSON-key row, row 12, keeps "RMTL" untouched).
$ws2.Range("F1").Value = "PMTL"
